# MarketBeatRank - 10th run: shift history columns right, add this week's
# columns (Jun_26 / Jun_27), note a new "Raises Target" event for Stifel
# Nicolaus, and append two newly-tracked analysts (Benchmark, Evercore ISI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room: insert three new columns before column B -----------------
# This pushes the existing Jun_17 / Jun_15 / Jun_13 / Jun_10 columns
# (B:E) to the right, landing on E:H.
$ws.Range("B:D").Insert()

# --- 2) New header row (B1:D1) for the two newest snapshots ------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3) Fill the new columns (rows 2-27) with "UN" like the older columns ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- 4) Highlight the Stifel Nicolaus (row 14) rating change in C14:D14 -----
# Match the same highlight fill already used on the other "rating changed"
# cells in this sheet (e.g. E10 / E17 before the column insert).
$ws.Range("C14:D14").Value = "6/18/2018,Raises Target,Hold -> Hold,$15.00 -> $18.00"
$ws.Range("C14:D14").Interior.Color = 13434828

# --- 5) Append two newly tracked analysts ------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
